$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.835.52"
$ws.Range("E2").Value = "  -0.97%  "
$ws.Range("D3").Value = "2.366.39"
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "'318.54"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.07%  "
$ws.Range("D6").Value = "'108.88"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.45%  "
$ws.Range("E7").Value = "  -2.26%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").Value = "'0.622"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.51%  "
$ws.Range("D10").Value = "'42.03"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.28%  "
$ws.Range("D11").Value = "'0.0928"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.43%  "
$ws.Range("D12").Value = "'8.62"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.38%  "
$ws.Range("D13").Value = "'1.01"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.29%  "
$ws.Range("D14").Value = "'0.107"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.07%  "
$ws.Range("D15").Value = "'16.18"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.50%  "
$ws.Range("D16").Value = "2.724.30"
$ws.Range("E16").Value = "  -1.20%  "
$ws.Range("D17").Value = "2.377.79"
$ws.Range("E17").Value = "  -1.29%  "
$ws.Range("D18").Value = "42.802.12"
$ws.Range("D19").Value = "'7.74"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.18%  "
$ws.Range("E20").Value = "  -1.70%  "
$ws.Range("D21").Value = "'76.31"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.93%  "
$ws.Range("D22").Value = "'3.68"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.69%  "
$ws.Range("D23").Value = "'257.17"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -6.08%  "
$ws.Range("D24").Value = "'2.33"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.37%  "
$ws.Range("D25").Value = "'9.46"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.76%  "
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("D27").Value = "'11.50"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.76%  "
$ws.Range("D28").Value = "'22.96"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.70%  "
$ws.Range("E29").Value = "  +2.40%  "
$ws.Range("D30").Value = "'37.48"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.08%  "
$ws.Range("D31").Value = "'172.41"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.15%  "
$ws.Range("D32").Value = "'0.0896"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.71%  "
$ws.Range("D33").Value = "'6.10"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.48%  "
$ws.Range("D34").Value = "'2.93"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -7.25%  "
$ws.Range("D35").Value = "'0.122"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +11.88%  "
$ws.Range("D36").Value = "'0.132"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.37%  "
$ws.Range("D37").Value = "'4.72"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.94%  "
$ws.Range("D38").Value = "'0.0365"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.00%  "
$ws.Range("D39").Value = "'3.93"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.48%  "
$ws.Range("D40").Value = "'2.70"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.48%  "
$ws.Range("E41").Value = "  +3.25%  "
$ws.Range("D42").Value = "'1.52"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.21%  "
$ws.Range("D43").Value = "'71.71"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.03%  "
$ws.Range("E44").Value = "  +0.01%  "
$ws.Range("D45").Value = "'12.43"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.82%  "
$ws.Range("D46").Value = "'5.61"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.67%  "
$ws.Range("D47").Value = "'112.97"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -8.44%  "
$ws.Range("E48").Value = "  -1.45%  "
$ws.Range("D49").Value = "'86.30"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.15%  "
$ws.Range("D50").Value = "'77.43"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.61%  "
$ws.Range("E51").Value = "  -1.18%  "
